$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Delete()

$ws.Range("A4").Clear()
$ws.Range("AJ4:AL4").Clear()

$ws.Range("AE3").Value = "SS"
$ws.Range("AF3").Value = "SS"
$ws.Range("AG3").Value = "SS"
$ws.Range("AH3").Value = "SS"
$ws.Range("AI3").Value = 20140101

$ws.Range("AE3").Select()
